$d = $word.ActiveDocument

# Change 1: update the URNFactor/ContextFactor line - drop "Provenance x" and
# pluralize PredicateFactor / ObjectFactor.
$d.Content.Find.Execute(
    "(URNFactor (ContextFactor x Provenance x SubjectKindFactors : class, SubjectFactors : instance, PredicateFactor : attribute, ObjectFactor : value))",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(URNFactor (ContextFactor x  SubjectKindFactors : class, SubjectFactors : instance, PredicateFactors : attribute, ObjectFactors : value))",
    2) | Out-Null

# Change 2: the old "Provenance: Entailment..." bullet becomes "Reified
# Statements:" and three new bullets (with the same list formatting) are
# inserted after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Provenance: Entailment. Provenance x Statement SPO Factors*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Provenance: Entailment...' bullet paragraph"
}

$target.Range.Text = "Reified Statements:"

$target.Range.InsertParagraphAfter() | Out-Null
$p2 = $target.Next()
$p2.Range.Text = "(StatementFactors : SPO Factors x Provenance (SubjectFactors x SubjectKind, PredicateFactors x PredicateKind, ObjectFactors x ObjectKind))"

$p2.Range.InsertParagraphAfter() | Out-Null
$p3 = $p2.Next()
$p3.Range.Text = "Provenance: Entailment. Provenance x reified SPO StatementFactors / source URNs. Statement entails / entailed by Factors product relations: transitive, reflexive, symmetrical (cause / effect, etc.)."

$p3.Range.InsertParagraphAfter() | Out-Null
$p4 = $p3.Next()
$p4.Range.Text = "Alignment / Matching: Factors of matching URNs / Statements aggregated by product."

Write-Output "done"
